$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 2 (UID=1) acts as the template for the new rows 3-17 (UID=2..16).
# Columns: A=UID, B=Source Plate Name, C=Source plate Type, D=Source Well (blank),
#          E=Destination Plate Name, F=Destination Plate Type, G=Destination Well,
#          H=Transfer Volume, I=Reagent
$sourcePlateName = $ws.Range("B2").Value()
$sourcePlateType = $ws.Range("C2").Value()
$destPlateName = $ws.Range("E2").Value()
$destPlateType = $ws.Range("F2").Value()
$transferVolume = $ws.Range("H2").Value()
$reagent = $ws.Range("I2").Value()

for ($i = 2; $i -le 16; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $sourcePlateName
    $ws.Cells.Item($row, 3).Value = $sourcePlateType
    $ws.Cells.Item($row, 5).Value = $destPlateName
    $ws.Cells.Item($row, 6).Value = $destPlateType
    $ws.Cells.Item($row, 7).Value = "A$i"
    $ws.Cells.Item($row, 8).Value = $transferVolume
    $ws.Cells.Item($row, 9).Value = $reagent
}
